function Get-ParagraphByText($doc, $pattern) {
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text -like $pattern) {
            return $para
        }
    }
    return $null
}

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr>'

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Phone Number: ... (optional)." paragraph -> split into 4 runs ending
#    in "(mandatory);" and add a new "Address ID" list paragraph after it.
# ---------------------------------------------------------------------------
$pPrPhone = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr>' + $rPr + '</w:pPr>'

$phone = Get-ParagraphByText $d "Phone Number:*"
$xmlPhone = '<w:p ' + $ns + '>' + $pPrPhone + `
    '<w:r>' + $rPr + '<w:t>Phone Number: the phone number of the user, used for notifications (</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>mandatory</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>)</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>;</w:t></w:r>' + `
    '</w:p>'
$phone.Range.InsertXML($xmlPhone)

$phone = Get-ParagraphByText $d "Phone Number:*"
$phone.Range.InsertParagraphAfter()

$d = $word.ActiveDocument
$phone = Get-ParagraphByText $d "Phone Number:*"
$newPara = $phone.Next()

$addrText = "Address ID: an unique identifier for the user" + [char]0x2019 + "s address."
$xmlAddressId = '<w:p ' + $ns + '>' + $pPrPhone + '<w:r>' + $rPr + '<w:t>' + $addrText + '</w:t></w:r></w:p>'
$newPara.Range.InsertXML($xmlAddressId)

# ---------------------------------------------------------------------------
# 2) "User ID: an unique identifier for the user that lives..." paragraph ->
#    rewritten to "Country: the name of the country where the user lives;"
#    then two new paragraphs (District, City) are inserted right after it,
#    before the "Street" paragraph.
# ---------------------------------------------------------------------------
$pPrAddr8 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr>' + $rPr + '</w:pPr>'

$userId = Get-ParagraphByText $d "User ID:*"
$xmlCountry = '<w:p ' + $ns + '>' + $pPrAddr8 + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Country: the name of the country </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>where the user lives</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>;</w:t></w:r>' + `
    '</w:p>'
$userId.Range.InsertXML($xmlCountry)

$country = Get-ParagraphByText $d "Country:*"
$country.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$country = Get-ParagraphByText $d "Country:*"
$districtPara = $country.Next()
$xmlDistrict = '<w:p ' + $ns + '>' + $pPrAddr8 + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">District: the name of the district/county </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>where the user lives</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>;</w:t></w:r>' + `
    '</w:p>'
$districtPara.Range.InsertXML($xmlDistrict)

$district = Get-ParagraphByText $d "District:*"
$district.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$district = Get-ParagraphByText $d "District:*"
$cityPara = $district.Next()
$xmlCity = '<w:p ' + $ns + '>' + $pPrAddr8 + `
    '<w:r>' + $rPr + '<w:t>City: the name of the city where the user lives;</w:t></w:r>' + `
    '</w:p>'
$cityPara.Range.InsertXML($xmlCity)

Write-Output "done"
